$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 136.45454
$ws.Range("I9").Value = 128.55556
$ws.Range("J9").Value = 172
$ws.Range("K9").Value = 128.55556
$ws.Range("L9").Value = 172
$ws.Range("M9").Value = 40.44443999999999
$ws.Range("N9").Value = -510

$ws.Range("H12").Value = 3300
$ws.Range("J12").Value = 200
$ws.Range("L12").Value = 200
$ws.Range("N12").Value = -540

$ws.Range("H113").Value = 41670148
$ws.Range("I113").Value = 58826204
$ws.Range("J113").Value = 5435.2856
$ws.Range("K113").Value = 58826204
$ws.Range("L113").Value = 5435.2856
$ws.Range("M113").Value = -58822950
$ws.Range("N113").Value = -11943.2856

$ws.Range("H129").Value = 257380
$ws.Range("J129").Value = 334528.06
$ws.Range("L129").Value = 1003584.18
$ws.Range("N129").Value = -1013584.18

$ws.Range("H132").Value = 2527.5557
$ws.Range("I132").Value = 2701.225
$ws.Range("K132").Value = 8103.674999999999
$ws.Range("M132").Value = -5573.674999999999

$ws.Range("H137").Value = 1227.5161
$ws.Range("I137").Value = 1292.2222
$ws.Range("K137").Value = 3876.6666
$ws.Range("M137").Value = -1326.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 659.6786
$ws.Range("I2").Value = 766.8946999999999
$ws.Range("K2").Value = 766.8946999999999
$ws.Range("M2").Value = -653.8946999999999

$ws.Range("H61").Value = 3160.8
$ws.Range("I61").Value = 2439
$ws.Range("K61").Value = 2439
$ws.Range("M61").Value = -2227

$ws.Range("H74").Value = 23810758
$ws.Range("I74").Value = 31250382
$ws.Range("K74").Value = 31250382
$ws.Range("M74").Value = -31249508

$ws.Range("H77").Value = 23810758
$ws.Range("I77").Value = 31250382
$ws.Range("K77").Value = 156251910
$ws.Range("M77").Value = -156247542

$ws.Range("H110").Value = 752.9091
$ws.Range("J110").Value = 1033
$ws.Range("L110").Value = 1033
$ws.Range("N110").Value = -5123

$ws.Range("H116").Value = 659.6786
$ws.Range("I116").Value = 766.8946999999999
$ws.Range("K116").Value = 766.8946999999999
$ws.Range("M116").Value = 1527.1053

$ws.Range("H136").Value = 3160.8
$ws.Range("I136").Value = 2439
$ws.Range("K136").Value = 7317
$ws.Range("M136").Value = -4767

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 659.6786
$ws.Range("I3").Value = 766.8946999999999
$ws.Range("K3").Value = 766.8946999999999
$ws.Range("M3").Value = -652.8946999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 711.4
$ws.Range("I10").Value = 689.25
$ws.Range("J10").Value = 800
$ws.Range("K10").Value = 689.25
$ws.Range("L10").Value = 800
$ws.Range("M10").Value = -550.25
$ws.Range("N10").Value = -1078

$ws.Range("H22").Value = 478.75
$ws.Range("I22").Value = 264.5
$ws.Range("K22").Value = 264.5
$ws.Range("M22").Value = 85.5

$ws.Range("H58").Value = 18061.3
$ws.Range("I58").Value = 1239.85
$ws.Range("J58").Value = 51704.2
$ws.Range("K58").Value = 1239.85
$ws.Range("L58").Value = 51704.2
$ws.Range("M58").Value = -1036.85
$ws.Range("N58").Value = -52110.2

$ws.Range("H107").Value = 1045.1482
$ws.Range("I107").Value = 375.875
$ws.Range("J107").Value = 2018.6364
$ws.Range("K107").Value = 375.875
$ws.Range("L107").Value = 2018.6364
$ws.Range("M107").Value = 1544.125
$ws.Range("N107").Value = -5858.6364

$ws.Range("H132").Value = 3483.9048
$ws.Range("I132").Value = 2618.1538
$ws.Range("J132").Value = 4890.75
$ws.Range("K132").Value = 7854.4614
$ws.Range("L132").Value = 14672.25
$ws.Range("M132").Value = -5324.4614
$ws.Range("N132").Value = -19732.25

$ws.Range("H134").Value = 1660.1765
$ws.Range("I134").Value = 1587.3572
$ws.Range("K134").Value = 4762.071599999999
$ws.Range("M134").Value = -2227.071599999999

$ws.Range("H136").Value = 18061.3
$ws.Range("I136").Value = 1239.85
$ws.Range("J136").Value = 51704.2
$ws.Range("K136").Value = 3719.55
$ws.Range("L136").Value = 155112.6
$ws.Range("M136").Value = -1169.55
$ws.Range("N136").Value = -160212.6

$ws.Range("H137").Value = 26638
$ws.Range("J137").Value = 26638
$ws.Range("L137").Value = 26638
$ws.Range("N137").Value = -36838

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1664.7778
$ws.Range("J5").Value = 2651.25
$ws.Range("L5").Value = 7953.75
$ws.Range("N5").Value = -8177.75

$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("N50").ClearContents()

$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("N53").ClearContents()

$ws.Range("H104").Value = 2450
$ws.Range("I104").Value = 2450
$ws.Range("K104").Value = 7350
$ws.Range("M104").Value = -4729

$ws.Range("H131").Value = 742.4400000000001
$ws.Range("J131").Value = 746.9091
$ws.Range("L131").Value = 2240.7273
$ws.Range("N131").Value = -12320.7273

$ws.Range("H132").Value = 695
$ws.Range("I132").Value = 695
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6255
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3725
$ws.Range("N132").ClearContents()

$ws.Range("H133").Value = 1129.6666
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H134").Value = 1432.68
$ws.Range("I134").Value = 1432.68
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4298.04
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 771.96
$ws.Range("N134").ClearContents()

$ws.Range("H135").Value = 1664.7778
$ws.Range("J135").Value = 2651.25
$ws.Range("L135").Value = 23861.25
$ws.Range("N135").Value = -28931.25

$ws.Range("H136").Value = 4043.5
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 4804.375
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 14413.125
$ws.Range("M136").Value = 2100
$ws.Range("N136").Value = -24613.125

$ws.Range("H137").Value = 25647262
$ws.Range("I137").Value = 1870
$ws.Range("J137").Value = 41675630
$ws.Range("K137").Value = 5610
$ws.Range("L137").Value = 125026890
$ws.Range("M137").Value = -510
$ws.Range("N137").Value = -125037090

$ws.Range("H141").Value = 2188.3333
$ws.Range("I141").Value = 2188.3333
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 6564.999899999999
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -1384.999899999999
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4153.8887
$ws.Range("I80").Value = 3260
$ws.Range("J80").Value = 4722.727
$ws.Range("K80").Value = 3260
$ws.Range("L80").Value = 4722.727
$ws.Range("M80").Value = -2262
$ws.Range("N80").Value = -6718.727

$ws.Range("H83").Value = 4153.8887
$ws.Range("I83").Value = 3260
$ws.Range("J83").Value = 4722.727
$ws.Range("K83").Value = 16300
$ws.Range("L83").Value = 23613.635
$ws.Range("M83").Value = -11308
$ws.Range("N83").Value = -33597.63499999999

$ws.Range("H122").Value = 60607840
$ws.Range("I122").Value = 19609026
$ws.Range("J122").Value = 200003800
$ws.Range("K122").Value = 58827078
$ws.Range("L122").Value = 600011400
$ws.Range("M122").Value = -58824628
$ws.Range("N122").Value = -600016300

$ws.Range("H125").Value = 20326
$ws.Range("J125").Value = 20326
$ws.Range("L125").Value = 20326
$ws.Range("N125").Value = -25246

$ws.Range("H126").Value = 3502.2273
$ws.Range("I126").Value = 2365.5625
$ws.Range("K126").Value = 7096.6875
$ws.Range("M126").Value = -4626.6875

$ws.Range("H132").Value = 41922.08
$ws.Range("I132").Value = 3314.4285
$ws.Range("K132").Value = 9943.2855
$ws.Range("M132").Value = -7413.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2888.5264
$ws.Range("I7").Value = 2481.353
$ws.Range("K7").Value = 2481.353
$ws.Range("M7").Value = -2369.353

$ws.Range("H126").Value = 2888.5264
$ws.Range("I126").Value = 2481.353
$ws.Range("K126").Value = 7444.059
$ws.Range("M126").Value = -4974.059

$ws.Range("H132").Value = 549953.0600000001
$ws.Range("I132").Value = 1205959.1
$ws.Range("J132").Value = 3281.4167
$ws.Range("K132").Value = 3617877.3
$ws.Range("L132").Value = 9844.250100000001
$ws.Range("M132").Value = -3615347.3
$ws.Range("N132").Value = -14904.2501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 48499.5
$ws.Range("J125").Value = 48499.5
$ws.Range("L125").Value = 48499.5
$ws.Range("N125").Value = -58339.5

$ws.Range("H126").Value = 1286.4667
$ws.Range("I126").Value = 1141.4584
$ws.Range("K126").Value = 3424.3752
$ws.Range("M126").Value = -954.3751999999999

$ws.Range("H132").Value = 1839.6111
$ws.Range("I132").Value = 1412.2222
$ws.Range("J132").Value = 2267
$ws.Range("K132").Value = 4236.6666
$ws.Range("L132").Value = 6801
$ws.Range("M132").Value = -1706.6666
$ws.Range("N132").Value = -11861
